# Updated cryptos list - apply latest price/volume values to worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.682.21"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.596.22"

$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").Value = "211.39"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("E8").Value = "  +0.31%  "

$ws.Range("E9").Value = "  +0.50%  "

$ws.Range("E10").Value = "  +0.06%  "

$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("D12").Value = "1.820.58"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "1.575.31"
$ws.Range("E13").Value = "  -1.66%  "

$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D16").Value = "65.10"
$ws.Range("E16").Value = "  +0.96%  "

$ws.Range("D17").Value = "26.648.52"
$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").Value = "0.0₃0762"
$ws.Range("E18").Value = "  +4.60%  "

$ws.Range("E19").Value = "  +0.34%  "

$ws.Range("D20").Value = "209.25"

$ws.Range("D21").Value = "7.05"
$ws.Range("E21").Value = "  +3.89%  "

$ws.Range("E22").Value = "  +0.84%  "

$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").Value = "8.92"
$ws.Range("E24").Value = "  +0.89%  "

$ws.Range("D25").Value = "142.85"

$ws.Range("E26").Value = "  +0.25%  "

$ws.Range("D27").Value = "7.10"
$ws.Range("E27").Value = "  -1.12%  "

$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("E29").Value = "  +0.68%  "

$ws.Range("D30").Value = "0.0517"
$ws.Range("E30").Value = "  +2.56%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("E33").Value = "  +1.01%  "

$ws.Range("D34").Value = "1.280.21"

$ws.Range("D35").Value = "0.616"
$ws.Range("E35").Value = "  -6.98%  "

$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("E39").Value = "  +18.00%  "

$ws.Range("E40").Value = "  -2.06%  "

$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  -0.68%  "

$ws.Range("E43").Value = "  -0.39%  "

$ws.Range("D44").Value = "62.87"
$ws.Range("E44").Value = "  -1.02%  "

$ws.Range("D45").Value = "1.733.01"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").Value = "90.77"
$ws.Range("E46").Value = "  +0.87%  "

$ws.Range("E47").Value = "  -2.53%  "

$ws.Range("D48").Value = "0.100"
$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("E50").Value = "  +0.31%  "

$ws.Range("D51").Value = "7.28"
$ws.Range("E51").Value = "  -2.12%  "
